# Minor refactoring to SpreadsheetImporter
#
# Core content change: the shared string "Machray" becomes "Machray Hall".
# It is referenced by every cell in column A that currently reads "Machray"
# (rows 9, 13, 14, 15, 16 on Sheet1) - update them all so the workbook ends
# up with a single shared string "Machray Hall" instead of "Machray".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value() -eq "Machray") {
        $cell.Value = "Machray Hall"
    }
}

# Cosmetic / view-state touch-ups that accompanied the rename in the
# original commit.
$ws.Tab.Color = 16777215
$ws.Rows.Item(24).RowHeight = 12.8
[void]$ws.Range("A18").Select()
